# Add "sigma_max" row values (row 8) across all columns B:H, and update
# the sheet view's scroll position / selection to match the saved state
# (topLeftCell B1, active selection H8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 ("sigma max") currently only has values in E8/F8 (=50 each).
# Fill in the remaining columns (B, C, D, G, H) with the same value (50)
# so the whole B8:H8 range is populated.
$ws.Range("B8").Value = 50
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 50
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 50
$ws.Range("G8").Value = 50
$ws.Range("H8").Value = 50

# Update the window scroll position (topLeftCell -> B1) and select H8,
# matching the sheetView/selection state recorded in the saved workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("H8").Select()
